$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-236) holds a date value (45175 -> 2023-09-06) that was
# bumped forward by two days (45177 -> 2023-09-08) for every data row.
$startRow = 2
$endRow = 236
$newValue = 45177

for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Cells.Item($r, 3).Value = $newValue
}
